# Update cryptocurrency price (D) and 1h-volume (E) data refreshed by the scraper
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.617.52"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.459.38"
$ws.Range("E3").Value = "  +0.98%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.29"
$ws.Range("E5").Value = "  +0.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.97"
$ws.Range("E6").Value = "  +1.68%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.511"
$ws.Range("E8").Value = "  +2.84%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.459.78"
$ws.Range("E9").Value = "  +1.23%  "
$ws.Range("E10").Value = "  +6.96%  "
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.329"
$ws.Range("E12").Value = "  -1.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.82"
$ws.Range("E13").Value = "  +1.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "68.523.98"
$ws.Range("E14").Value = "  +0.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000170"
$ws.Range("E15").Value = "  +2.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.47"
$ws.Range("E16").Value = "  +2.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "10.54"
$ws.Range("E17").Value = "  -1.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "338.06"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.90"
$ws.Range("E19").Value = "  -1.96%  "
$ws.Range("E20").Value = "  +2.53%  "
$ws.Range("E21").Value = "  +3.93%  "
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.61"
$ws.Range("E23").Value = "  +0.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.66"
$ws.Range("E24").Value = "  +1.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.17"
$ws.Range("E25").Value = "  +2.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0₃0820"
$ws.Range("E26").Value = "  +0.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.21"
$ws.Range("E27").Value = "  +2.47%  "
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "424.78"
$ws.Range("E29").Value = "  +0.96%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.14"
$ws.Range("E30").Value = "  +0.79%  "
$ws.Range("E31").Value = "  +0.99%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "159.49"
$ws.Range("E32").Value = "  +1.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.96"
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.106"
$ws.Range("E35").Value = "  -1.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.83"
$ws.Range("E36").Value = "  +1.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.40"
$ws.Range("E37").Value = "  +2.13%  "
$ws.Range("E38").Value = "  -0.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.47"
$ws.Range("E39").Value = "  +0.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.06"
$ws.Range("E40").Value = "  -0.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.04"
$ws.Range("E41").Value = "  +0.69%  "
$ws.Range("E42").Value = "  +2.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "129.99"
$ws.Range("E43").Value = "  -2.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0721"
$ws.Range("E44").Value = "  +1.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.482"
$ws.Range("E45").Value = "  +1.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.562"
$ws.Range("E46").Value = "  +1.25%  "
$ws.Range("E47").Value = "  +2.20%  "
$ws.Range("E48").Value = "  +0.93%  "
$ws.Range("E49").Value = "  -0.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.94"
$ws.Range("E50").Value = "  -3.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.75"
$ws.Range("E51").Value = "  -0.41%  "
